$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the changed Price (D) and Volume% (E) cells keep their original text (string) type
# instead of being auto-converted to numbers/percentages by Excel.
$cells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","E17","E18","E19","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","E27","D39","E39","D40","E40","D41","E41","E42","D43","E43","D44","E44","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "303.45"
$ws.Range("E2").Value = "0.67%"
$ws.Range("D3").Value = "35.57"
$ws.Range("E3").Value = "11.14%"
$ws.Range("D4").Value = "5.083"
$ws.Range("E4").Value = "1.47%"
$ws.Range("D5").Value = "0.07775"
$ws.Range("E5").Value = "-0.45%"
$ws.Range("D6").Value = "2.255"
$ws.Range("E6").Value = "-3.30%"
$ws.Range("D7").Value = "8.064"
$ws.Range("E7").Value = "1.06%"
$ws.Range("D8").Value = "4.047"
$ws.Range("E8").Value = "4.09%"
$ws.Range("D9").Value = "0.9288"
$ws.Range("E9").Value = "-0.63%"
$ws.Range("D10").Value = "0.09320"
$ws.Range("E10").Value = "-8.63%"
$ws.Range("D11").Value = "0.1833"
$ws.Range("E11").Value = "2.67%"
$ws.Range("D12").Value = "0.08576"
$ws.Range("E12").Value = "0.60%"
$ws.Range("D13").Value = "0.03721"
$ws.Range("E13").Value = "12.62%"
$ws.Range("D14").Value = "0.09930"
$ws.Range("E14").Value = "0.39%"
$ws.Range("D15").Value = "0.001476"
$ws.Range("E15").Value = "0.19%"
$ws.Range("D16").Value = "0.005750"
$ws.Range("E16").Value = "0.01%"
$ws.Range("E17").Value = "-0.23%"
$ws.Range("E18").Value = "0.06%"
$ws.Range("E19").Value = "3.00%"
$ws.Range("E20").Value = "-1.47%"
$ws.Range("D21").Value = "4.551"
$ws.Range("E21").Value = "5.15%"
$ws.Range("D22").Value = "0.2237"
$ws.Range("E22").Value = "7.48%"
$ws.Range("D23").Value = "0.04672"
$ws.Range("E23").Value = "1.49%"
$ws.Range("D24").Value = "0.001233"
$ws.Range("E24").Value = "1.16%"
$ws.Range("D25").Value = "0.004531"
$ws.Range("E25").Value = "3.06%"
$ws.Range("D26").Value = "0.0001303"
$ws.Range("E26").Value = "0.63%"
$ws.Range("E27").Value = "-19.87%"
$ws.Range("D39").Value = "0.01772"
$ws.Range("E39").Value = "3.54%"
$ws.Range("D40").Value = "0.04728"
$ws.Range("E40").Value = "-1.24%"
$ws.Range("D41").Value = "0.007969"
$ws.Range("E41").Value = "2.67%"
$ws.Range("E42").Value = "0.56%"
$ws.Range("D43").Value = "0.007778"
$ws.Range("E43").Value = "-20.34%"
$ws.Range("D44").Value = "0.002226"
$ws.Range("E44").Value = "7.13%"
$ws.Range("E45").Value = "-5.67%"
$ws.Range("D46").Value = "0.00006204"
$ws.Range("E46").Value = "1.98%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.62%"
$ws.Range("D48").Value = "5.353"
$ws.Range("E48").Value = "91.61%"
$ws.Range("D49").Value = "0.002696"
$ws.Range("E49").Value = "35.77%"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").Value = "0.62%"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").Value = "0.62%"
